$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update B3 value (PSC) from 0 to 47 -> this recalculates F3 formula to 999
$ws.Range("B3").Value = 47

# Update the active cell selection to B4
$ws.Range("B4").Select()
